$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K" = strikeouts) is being regenerated to reflect actual
# strikeout counts instead of the previous "Strike#" pitch-count values.
$kValues = @{
    2  = 6
    3  = 7
    4  = 9
    5  = 3
    6  = 11
    7  = 3
    8  = 6
    9  = 5
    10 = 12
    11 = 11
    12 = 6
    13 = 3
    14 = 2
    15 = 8
    16 = 11
    17 = 12
    18 = 5
    19 = 8
    20 = 5
    21 = 5
    22 = 5
    23 = 10
    24 = 2
    25 = 8
    26 = 12
    27 = 9
    28 = 9
    29 = 6
    30 = 7
    31 = 6
    32 = 4
    33 = 3
    34 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
